$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Adhoc")

# Seed row 3 from row 2 (same text/number cell types, no coercion or stray
# number-format styles), then overwrite the cells that actually differ.
$ws.Range("A2:F2").Copy()
$ws.Range("A3").PasteSpecial()

$ws.Range("A3").Value = "M-004"
$ws.Range("C3").Value = 75
$ws.Range("D3").Value = "paid"
$ws.Range("E3").Value = "Reimbursement"
$ws.Range("F3").Value = "Sample paid adhoc payment"
